# Estadisticos Matutinos 15 Oct
# Updates the statistics rows on the three "Estadisticos" sheets and adds a
# rescue-exam ("Rescatables") record for one more student.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P: Reprobados/Aprobados/Por_Apro/Promedio for row 2 ---
$ws1P = $wb.Worksheets.Item("Estadisticos 1P")
$ws1P.Range("D2").Value = 14
$ws1P.Range("F2").Value = 19
$ws1P.Range("G2").Value = 57.58
$ws1P.Range("H2").Value = 7.4

# --- Estadisticos 2P: only Reprobados (E2) updates ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("E2").Value = 19

# --- Estadisticos Final: same update pattern as 1P ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("D2").Value = 14
$wsFinal.Range("F2").Value = 19
$wsFinal.Range("G2").Value = 57.58
$wsFinal.Range("H2").Value = 7.4

# --- Rescatables: add the new student row ---
$wsResc = $wb.Worksheets.Item("Rescatables")
$wsResc.Range("A2").Value = 20330051920072
$wsResc.Range("B2").Value = "CARRERA"
$wsResc.Range("C2").Value = "CASTAÑEDA"
$wsResc.Range("D2").Value = "PAUL ARAVIER"
$wsResc.Range("E2").Value = "MANTIENE LOS GENERADORES DE CA Y CC"
$wsResc.Range("F2").Value = "3AEV"
$wsResc.Range("G2").Value = 6
